$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: financial periods (header labels) - shift left, append new period at the end
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish dates - shift left, append new date at the end
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-02-10 (9)"
$ws.Range("F9").Value = "1401-04-28 (3)"
$ws.Range("G9").Value = "1401-08-30 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-12 (10)"
$ws.Range("J9").Value = "1401-04-28 (2)"
$ws.Range("K9").Value = "1401-08-30 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-12 (2)"

# Rows 11-27: financial data values (updated per new database + read_price algorithm change)
# Row 11
$ws.Range("D11").Value = 16186
$ws.Range("E11").Value = 20107
$ws.Range("F11").Value = 7173
$ws.Range("G11").Value = 16771
$ws.Range("H11").Value = 25021
$ws.Range("I11").Value = 30949
$ws.Range("J11").Value = 10498
$ws.Range("K11").Value = 21171
$ws.Range("L11").Value = 30537
$ws.Range("M11").Value = 34419

# Row 12
$ws.Range("D12").Value = -10223
$ws.Range("E12").Value = -13228
$ws.Range("F12").Value = -4288
$ws.Range("G12").Value = -8807
$ws.Range("H12").Value = -13646
$ws.Range("I12").Value = -18884
$ws.Range("J12").Value = -5518
$ws.Range("K12").Value = -11214
$ws.Range("L12").Value = -17123
$ws.Range("M12").Value = -21362

# Row 13
$ws.Range("D13").Value = 5963
$ws.Range("E13").Value = 6879
$ws.Range("F13").Value = 2885
$ws.Range("G13").Value = 7963
$ws.Range("H13").Value = 11375
$ws.Range("I13").Value = 12065
$ws.Range("J13").Value = 4981
$ws.Range("K13").Value = 9957
$ws.Range("L13").Value = 13414
$ws.Range("M13").Value = 13057

# Row 14
$ws.Range("D14").Value = -1137
$ws.Range("E14").Value = -1535
$ws.Range("F14").Value = -398
$ws.Range("G14").Value = -794
$ws.Range("H14").Value = -1058
$ws.Range("I14").Value = -1664
$ws.Range("J14").Value = -474
$ws.Range("K14").Value = -1188
$ws.Range("L14").Value = -1557
$ws.Range("M14").Value = -2108

# Row 15
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"
$ws.Range("I15").Value = "-"
$ws.Range("J15").Value = "-"
$ws.Range("K15").Value = "-"
$ws.Range("L15").Value = "-"
$ws.Range("M15").Value = "-"

# Row 16
$ws.Range("D16").Value = 366
$ws.Range("E16").Value = 363
$ws.Range("F16").Value = -7
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -31
$ws.Range("I16").Value = -1305
$ws.Range("J16").Value = -35
$ws.Range("K16").Value = 106
$ws.Range("L16").Value = -113
$ws.Range("M16").Value = 323

# Row 17
$ws.Range("D17").Value = 5192
$ws.Range("E17").Value = 5707
$ws.Range("F17").Value = 2479
$ws.Range("G17").Value = 7179
$ws.Range("H17").Value = 10286
$ws.Range("I17").Value = 9096
$ws.Range("J17").Value = 4472
$ws.Range("K17").Value = 8874
$ws.Range("L17").Value = 11743
$ws.Range("M17").Value = 11271

# Row 18
$ws.Range("D18").Value = -68
$ws.Range("E18").Value = -66
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = "-"
$ws.Range("I18").Value = "-"
$ws.Range("J18").Value = "-"
$ws.Range("K18").Value = "-"
$ws.Range("L18").Value = "-"
$ws.Range("M18").Value = -21

# Row 19
$ws.Range("D19").Value = 1638
$ws.Range("E19").Value = 1565
$ws.Range("F19").Value = 169
$ws.Range("G19").Value = 438
$ws.Range("H19").Value = 628
$ws.Range("I19").Value = 900
$ws.Range("J19").Value = 575
$ws.Range("K19").Value = 793
$ws.Range("L19").Value = 1206
$ws.Range("M19").Value = 1395

# Row 20
$ws.Range("D20").Value = 6762
$ws.Range("E20").Value = 7206
$ws.Range("F20").Value = 2648
$ws.Range("G20").Value = 7617
$ws.Range("H20").Value = 10914
$ws.Range("I20").Value = 9996
$ws.Range("J20").Value = 5047
$ws.Range("K20").Value = 9668
$ws.Range("L20").Value = 12950
$ws.Range("M20").Value = 12646

# Row 21
$ws.Range("D21").Value = -1189
$ws.Range("E21").Value = -725
$ws.Range("F21").Value = -490
$ws.Range("G21").Value = -1432
$ws.Range("H21").Value = -2056
$ws.Range("I21").Value = -1293
$ws.Range("J21").Value = -864
$ws.Range("K21").Value = -1240
$ws.Range("L21").Value = -1963
$ws.Range("M21").Value = -1160

# Row 22
$ws.Range("D22").Value = 5573
$ws.Range("E22").Value = 6481
$ws.Range("F22").Value = 2158
$ws.Range("G22").Value = 6186
$ws.Range("H22").Value = 8858
$ws.Range("I22").Value = 8703
$ws.Range("J22").Value = 4183
$ws.Range("K22").Value = 8428
$ws.Range("L22").Value = 10986
$ws.Range("M22").Value = 11486

# Row 23
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"
$ws.Range("I23").Value = "-"
$ws.Range("J23").Value = "-"
$ws.Range("K23").Value = "-"
$ws.Range("L23").Value = "-"
$ws.Range("M23").Value = "-"

# Row 24
$ws.Range("D24").Value = 5573
$ws.Range("E24").Value = 6481
$ws.Range("F24").Value = 2158
$ws.Range("G24").Value = 6186
$ws.Range("H24").Value = 8858
$ws.Range("I24").Value = 8703
$ws.Range("J24").Value = 4183
$ws.Range("K24").Value = 8428
$ws.Range("L24").Value = 10986
$ws.Range("M24").Value = 11486

# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0

# Row 26
$ws.Range("D26").Value = 5001
$ws.Range("E26").Value = 4865
$ws.Range("F26").Value = 4715
$ws.Range("G26").Value = 4439
$ws.Range("H26").Value = 4228
$ws.Range("I26").Value = 4169
$ws.Range("J26").Value = 3743
$ws.Range("K26").Value = 3643
$ws.Range("L26").Value = 3459
$ws.Range("M26").Value = 3117

# Row 27
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
